$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25 (A=0..23), columns B,C,D,F,G,I,L,N
# (columns E,H,J,K,M,O remain 0 and are left untouched)
$cols = @("B","C","D","F","G","I","L","N")

$data = @(
    @(20.75849094258156, 10.12805355882275, 6.467920147422119, 35.40359363320086, 3.672526356768738, 27.92570568027111, 10.75703469696627, 18.63971877444789),
    @(20.22724841127363, 9.458311091792739, 6.497985263345234, 35.1294210752789, 3.676593035323175, 27.91394961070971, 10.73514770621622, 18.71375881098562),
    @(19.90006335050689, 9.021638504189058, 6.517380437600713, 34.97206906234909, 3.679218369640668, 27.91445564673154, 10.72405592769416, 18.76120133611924),
    @(19.76669971709113, 8.83727200409861, 6.525519491971831, 34.91076418560796, 3.680320623928991, 27.91659832118633, 10.72012839855745, 18.78103430432746),
    @(19.74455932904738, 8.806269830280259, 6.526885198250223, 34.90075616327469, 3.680505613600339, 27.91707086846773, 10.71951207578068, 18.78435778026279),
    @(19.89826460137211, 9.019178088826694, 6.517489250473531, 34.97123081003471, 3.679233103637027, 27.91447671195526, 10.72400055815159, 18.76146678513027),
    @(20.57565445837599, 9.902396068041282, 6.478092662558245, 35.30681171235497, 3.673901979905101, 27.92004588680834, 10.74900249472341, 18.66483736883295),
    @(21.8869491207778, 11.43337511996887, 6.408244818876832, 36.0494566194727, 3.664460449856211, 27.99245774308825, 10.81653002123228, 18.49100282987307),
    @(22.82849582077362, 12.43692161156065, 6.361429349301098, 36.64273515748153, 3.658133019832367, 28.08333162352291, 10.87722414233331, 18.37273307190751),
    @(23.24990386895447, 12.86740471487808, 6.341106546144506, 36.92205674469118, 3.655385057722482, 28.13286120245172, 10.9071895579083, 18.32095953344001),
    @(23.4083233500663, 13.02669516894596, 6.333550734874353, 37.02910251857443, 3.654363095276061, 28.15279250621844, 10.91887013779159, 18.30164435699067),
    @(23.3742589010078, 12.99255443914852, 6.335171788063436, 37.00599305822108, 3.654582366497933, 28.14844771433924, 10.9163397816342, 18.3057913330996),
    @(23.26296125505872, 12.88058416133245, 6.340482119444252, 36.93083844872226, 3.655300607562574, 28.13447742211881, 10.90814387830964, 18.31936465069797),
    @(23.19463266797792, 12.81151472208316, 6.343753081379272, 36.8849672859664, 3.655742973645043, 28.12607320872588, 10.90316689246076, 18.32771647264815),
    @(22.80080238765182, 12.40826631553049, 6.362777083665111, 36.62466358490032, 3.658315219454664, 28.08025940676848, 10.87531277960348, 18.37615729243036),
    @(22.55730938006644, 12.15423223833088, 6.374697068990559, 36.46733356714605, 3.659926524572362, 28.0542508176034, 10.85882477838623, 18.40639266067418),
    @(22.4166176315878, 12.00566681955176, 6.381644787421051, 36.37773492200435, 3.660865587491614, 28.04006253952751, 10.84956336583013, 18.42397420219207),
    @(22.36887699217103, 11.95494343378225, 6.384012909408507, 36.34755432361965, 3.661185651392647, 28.03539112546769, 10.84646590366549, 18.42995984643348),
    @(22.58329721018989, 12.18152820226887, 6.373418680288943, 36.48398975143267, 3.659753728055677, 28.0569396741379, 10.86055701035524, 18.40315430008492),
    @(23.29568476765649, 12.91357346184258, 6.338918547169513, 36.95287932122488, 3.655089138083226, 28.1385489580279, 10.910542212523, 18.3153699686898),
    @(23.75445073481101, 13.37030925166524, 6.317186706186324, 37.26670354185335, 3.652149098870814, 28.19873505019195, 10.94515049809392, 18.25968956986924),
    @(23.51027397169618, 13.1285196461028, 6.328710740372431, 37.09856304601414, 3.653708362168267, 28.16598700842701, 10.92650378384164, 18.28925287082201),
    @(22.57155029345879, 12.16919554098419, 6.373996344661339, 36.47645682888835, 3.659831809787299, 28.05572166099042, 10.85977318961016, 18.40461774351326),
    @(21.53526547527684, 11.0405925837945, 6.426348976573165, 35.83988112224133, 3.666907053035307, 27.96626011583165, 10.79629976739141, 18.53636379482689)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $data[$i][$j]
    }
}
